$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (record #24) - mirrors row 2 (Ryanair FR3693 from Birmingham) but for Saturday, Jan 14
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(25, 3).Value = "9:20 AM"
$ws.Cells.Item(25, 4).Value = "FR3693"
$ws.Cells.Item(25, 5).Value = "Birmingham"
$ws.Cells.Item(25, 6).Value = "(BHX)"
$ws.Cells.Item(25, 7).Value = "Ryanair "
$ws.Cells.Item(25, 8).Value = "B738"
$ws.Cells.Item(25, 9).Value = "(EI-EVH)"
$ws.Cells.Item(25, 10).Value = "8:55 AM"
$ws.Cells.Item(25, 11).Borders.LineStyle = -4142
$ws.Cells.Item(25, 12).Value = "0 hours, -25 minutes"
$ws.Cells.Item(25, 13).Borders.LineStyle = -4142

# Row 26 (record #25) - mirrors row 3 (Ryanair FR2473 from London) but for Saturday, Jan 14
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(26, 3).Value = "9:30 AM"
$ws.Cells.Item(26, 4).Value = "FR2473"
$ws.Cells.Item(26, 5).Value = "London"
$ws.Cells.Item(26, 6).Value = "(STN)"
$ws.Cells.Item(26, 7).Value = "Ryanair "
$ws.Cells.Item(26, 8).Value = "B38M"
$ws.Cells.Item(26, 9).Value = "(EI-HMS)"
$ws.Cells.Item(26, 10).Value = "9:13 AM"
$ws.Cells.Item(26, 11).Borders.LineStyle = -4142
$ws.Cells.Item(26, 12).Value = "0 hours, -17 minutes"
$ws.Cells.Item(26, 13).Borders.LineStyle = -4142
